$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the new "discDateWithTimestamp" column (AL) header first so the new
#    shared string lands in the same order as the target file.
$ws.Range("AL1").Value = "discDateWithTimestamp"
# Match the other header cells' style (row default / s="2").
$ws.Range("AH1").Copy()
$ws.Range("AL1").PasteSpecial(-4122) | Out-Null
$ws.Range("AL1").Value = "discDateWithTimestamp"

# 2. Update the discDate value (T2) to the new date - stays a text cell.
$ws.Range("T2").Value = "05-12-2021"

# 3. Populate the new discDateWithTimestamp data cell (AL2) with the ISO
#    timestamp and copy T2's number format (Text "@") onto it.
$ws.Range("AL2").Value = "2021-05-12T00:00:00"
$ws.Range("T2").Copy()
$ws.Range("AL2").PasteSpecial(-4122) | Out-Null
$ws.Range("AL2").Value = "2021-05-12T00:00:00"

# 4. Re-format the two zip code cells to an integer number format, left
#    aligned (mirrors the new cellXfs entries 20/21 in the diff).
$ws.Range("S2").NumberFormat = "0"
$ws.Range("AE2").NumberFormat = "0"
$ws.Range("AE2").HorizontalAlignment = -4131

# 5. Size the new column and move the frozen-pane selection to the new cell.
$ws.Columns.Item(38).ColumnWidth = 20
$ws.Range("AL2").Select() | Out-Null
